$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking "Price" text values stay text (column D is stored as text in the source data)
$textPriceCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D13", "D14", "D15", "D16", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values scraped for this run
$ws.Range("D2").Value = "26.330.23"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.691.58"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "218.69"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "0.5274"
$ws.Range("E6").Value = "  +4.37%  "
$ws.Range("D7").Value = "1.008"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.2708"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").Value = "0.06433"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("D10").Value = "22.09"
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").Value = "1.718.88"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").Value = "4.569"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").Value = "0.5862"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "0.000008520"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "64.61"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "26.359.52"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "10.91"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "189.98"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").Value = "6.224"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "144.78"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "7.704"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").Value = "0.1233"
$ws.Range("E26").Value = "  +5.31%  "
$ws.Range("D27").Value = "15.89"
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("D28").Value = "0.06705"
$ws.Range("E28").Value = "  +15.75%  "
$ws.Range("D29").Value = "1.360"
$ws.Range("E29").Value = "  +6.34%  "
$ws.Range("D30").Value = "1.331"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").Value = "3.593"
$ws.Range("E31").Value = "  +2.34%  "
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("D35").Value = "0.6246"
$ws.Range("E35").Value = "  +4.38%  "
$ws.Range("D36").Value = "2.395"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").Value = "2.703"
$ws.Range("E37").Value = "  +2.32%  "
$ws.Range("D38").Value = "6.370"
$ws.Range("E38").Value = "  +5.96%  "
$ws.Range("D39").Value = "1.114.40"
$ws.Range("E39").Value = "  +3.87%  "
$ws.Range("D40").Value = "0.01623"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").Value = "0.8885"
$ws.Range("E41").Value = "  +3.31%  "
$ws.Range("D42").Value = "1.018"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").Value = "100.94"
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("D44").Value = "1.838.78"
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").Value = "0.00000000116"
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("D46").Value = "57.01"
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("D47").Value = "1.012"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "8.175"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").Value = "0.05271"
$ws.Range("D50").Value = "0.4304"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "6.078"
$ws.Range("E51").Value = "  +3.87%  "
